$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.365760326385498
$ws.Range("B1").Value = 2.620818138122559
$ws.Range("C1").Value = 2.811779260635376
$ws.Range("D1").Value = 3.276329755783081
$ws.Range("E1").Value = 0.8088361620903015
